$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A128").Value = 1004
$ws.Range("B128").Value = 279
$ws.Range("C128").Value = "V"
$ws.Range("D128").Value = "LasVegas"
$ws.Range("E128").Value = 0
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 14
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 14
$ws.Range("J128").Value = 52.5
$ws.Range("K128").Value = 51.5
$ws.Range("L128").Value = 150
$ws.Range("M128").Value = 2.5

$ws.Range("A129").Value = 1004
$ws.Range("B129").Value = 280
$ws.Range("C129").Value = "H"
$ws.Range("D129").Value = "LAChargers"
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 14
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 7
$ws.Range("I129").Value = 28
$ws.Range("J129").Value = 3
$ws.Range("K129").Value = 3
$ws.Range("L129").Value = -170
$ws.Range("M129").Value = 24

$ws.Range("A130").Value = 1007
$ws.Range("B130").Value = 301
$ws.Range("C130").Value = "V"
$ws.Range("D130").Value = "LARams"
$ws.Range("E130").Value = 0
$ws.Range("F130").Value = 3
$ws.Range("G130").Value = 13
$ws.Range("H130").Value = 10
$ws.Range("I130").Value = 26
$ws.Range("J130").Value = 52.5
$ws.Range("K130").Value = 3
$ws.Range("L130").Value = -140
$ws.Range("M130").Value = 1

$ws.Range("A131").Value = 1007
$ws.Range("B131").Value = 302
$ws.Range("C131").Value = "H"
$ws.Range("D131").Value = "Seattle"
$ws.Range("E131").Value = 0
$ws.Range("F131").Value = 7
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 10
$ws.Range("I131").Value = 17
$ws.Range("J131").Value = "pk"
$ws.Range("K131").Value = 53.5
$ws.Range("L131").Value = 120
$ws.Range("M131").Value = 26.5

$ws.Range("A132").Value = 1010
$ws.Range("B132").Value = 451
$ws.Range("C132").Value = "V"
$ws.Range("D132").Value = "NYJets"
$ws.Range("E132").Value = 0
$ws.Range("F132").Value = 3
$ws.Range("G132").Value = 6
$ws.Range("H132").Value = 11
$ws.Range("I132").Value = 20
$ws.Range("J132").Value = 44.5
$ws.Range("K132").Value = 45.5
$ws.Range("L132").Value = 130
$ws.Range("M132").Value = 0.5

$ws.Range("A133").Value = 1010
$ws.Range("B133").Value = 452
$ws.Range("C133").Value = "H"
$ws.Range("D133").Value = "Atlanta"
$ws.Range("E133").Value = 10
$ws.Range("F133").Value = 10
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 7
$ws.Range("I133").Value = 27
$ws.Range("J133").Value = 3.5
$ws.Range("K133").Value = 2.5
$ws.Range("L133").Value = -150
$ws.Range("M133").Value = 22.5

$ws.Range("A134").Value = 1010
$ws.Range("B134").Value = 453
$ws.Range("C134").Value = "V"
$ws.Range("D134").Value = "Miami"
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 7
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 17
$ws.Range("J134").Value = 50
$ws.Range("K134").Value = 48
$ws.Range("L134").Value = 450
$ws.Range("M134").Value = 23.5

$ws.Range("A135").Value = 1010
$ws.Range("B135").Value = 454
$ws.Range("C135").Value = "H"
$ws.Range("D135").Value = "TampaBay"
$ws.Range("E135").Value = 7
$ws.Range("F135").Value = 17
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 21
$ws.Range("I135").Value = 45
$ws.Range("J135").Value = 9
$ws.Range("K135").Value = 12
$ws.Range("L135").Value = -600
$ws.Range("M135").Value = 4.5

$ws.Range("A136").Value = 1010
$ws.Range("B136").Value = 455
$ws.Range("C136").Value = "V"
$ws.Range("D136").Value = "Philadelphia"
$ws.Range("E136").Value = 3
$ws.Range("F136").Value = 3
$ws.Range("G136").Value = 7
$ws.Range("H136").Value = 8
$ws.Range("I136").Value = 21
$ws.Range("J136").Value = 47
$ws.Range("K136").Value = 46.5
$ws.Range("L136").Value = 120
$ws.Range("M136").Value = 0.5

$ws.Range("A137").Value = 1010
$ws.Range("B137").Value = 456
$ws.Range("C137").Value = "H"
$ws.Range("D137").Value = "Carolina"
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = 5
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 3
$ws.Range("I137").Value = 18
$ws.Range("J137").Value = 3.5
$ws.Range("K137").Value = 2.5
$ws.Range("L137").Value = -140
$ws.Range("M137").Value = 21.5

$ws.Range("A138").Value = 1010
$ws.Range("B138").Value = 457
$ws.Range("C138").Value = "V"
$ws.Range("D138").Value = "NewOrleans"
$ws.Range("E138").Value = 7
$ws.Range("F138").Value = 13
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 13
$ws.Range("I138").Value = 33
$ws.Range("J138").Value = 1.5
$ws.Range("K138").Value = 2.5
$ws.Range("L138").Value = -135
$ws.Range("M138").Value = 21.5

$ws.Range("A139").Value = 1010
$ws.Range("B139").Value = 458
$ws.Range("C139").Value = "H"
$ws.Range("D139").Value = "Washington"
$ws.Range("E139").Value = 6
$ws.Range("F139").Value = 7
$ws.Range("G139").Value = 3
$ws.Range("H139").Value = 6
$ws.Range("I139").Value = 22
$ws.Range("J139").Value = 43.5
$ws.Range("K139").Value = 43
$ws.Range("L139").Value = 115
$ws.Range("M139").Value = 1.5

$ws.Range("A140").Value = 1010
$ws.Range("B140").Value = 459
$ws.Range("C140").Value = "V"
$ws.Range("D140").Value = "Tennessee"
$ws.Range("E140").Value = 7
$ws.Range("F140").Value = 17
$ws.Range("G140").Value = 7
$ws.Range("H140").Value = 6
$ws.Range("I140").Value = 37
$ws.Range("J140").Value = 7.5
$ws.Range("K140").Value = 4
$ws.Range("L140").Value = -200
$ws.Range("M140").Value = 24

$ws.Range("A141").Value = 1010
$ws.Range("B141").Value = 460
$ws.Range("C141").Value = "H"
$ws.Range("D141").Value = "Jacksonville"
$ws.Range("E141").Value = 6
$ws.Range("F141").Value = 7
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 6
$ws.Range("I141").Value = 19
$ws.Range("J141").Value = 51.5
$ws.Range("K141").Value = 48.5
$ws.Range("L141").Value = 175
$ws.Range("M141").Value = 0.5

$ws.Range("A142").Value = 1010
$ws.Range("B142").Value = 461
$ws.Range("C142").Value = "V"
$ws.Range("D142").Value = "Detroit"
$ws.Range("E142").Value = 3
$ws.Range("F142").Value = 3
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 11
$ws.Range("I142").Value = 17
$ws.Range("J142").Value = 49.5
$ws.Range("K142").Value = 49.5
$ws.Range("L142").Value = 350
$ws.Range("M142").Value = 23.5

$ws.Range("A143").Value = 1010
$ws.Range("B143").Value = 462
$ws.Range("C143").Value = "H"
$ws.Range("D143").Value = "Minnesota"
$ws.Range("E143").Value = 3
$ws.Range("F143").Value = 10
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 6
$ws.Range("I143").Value = 19
$ws.Range("J143").Value = 8
$ws.Range("K143").Value = 10
$ws.Range("L143").Value = -420
$ws.Range("M143").Value = 4.5

$ws.Range("A144").Value = 1010
$ws.Range("B144").Value = 463
$ws.Range("C144").Value = "V"
$ws.Range("D144").Value = "Denver"
$ws.Range("E144").Value = 3
$ws.Range("F144").Value = 3
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 13
$ws.Range("I144").Value = 19
$ws.Range("J144").Value = 1.5
$ws.Range("K144").Value = 1
$ws.Range("L144").Value = -125
$ws.Range("M144").Value = 3

$ws.Range("A145").Value = 1010
$ws.Range("B145").Value = 464
$ws.Range("C145").Value = "H"
$ws.Range("D145").Value = "Pittsburgh"
$ws.Range("E145").Value = 7
$ws.Range("F145").Value = 10
$ws.Range("G145").Value = 7
$ws.Range("H145").Value = 3
$ws.Range("I145").Value = 27
$ws.Range("J145").Value = 42
$ws.Range("K145").Value = 40
$ws.Range("L145").Value = 105
$ws.Range("M145").Value = 20

$ws.Range("A146").Value = 1010
$ws.Range("B146").Value = 465
$ws.Range("C146").Value = "V"
$ws.Range("D146").Value = "GreenBay"
$ws.Range("E146").Value = 0
$ws.Range("F146").Value = 16
$ws.Range("G146").Value = 3
$ws.Range("H146").Value = 3
$ws.Range("I146").Value = 25
$ws.Range("J146").Value = 3
$ws.Range("K146").Value = 1
$ws.Range("L146").Value = -130
$ws.Range("M146").Value = 1

$ws.Range("A147").Value = 1010
$ws.Range("B147").Value = 466
$ws.Range("C147").Value = "H"
$ws.Range("D147").Value = "Cincinnati"
$ws.Range("E147").Value = 7
$ws.Range("F147").Value = 7
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 8
$ws.Range("I147").Value = 22
$ws.Range("J147").Value = 48.5
$ws.Range("K147").Value = 50
$ws.Range("L147").Value = 110
$ws.Range("M147").Value = 24.5

$ws.Range("A148").Value = 1010
$ws.Range("B148").Value = 467
$ws.Range("C148").Value = "V"
$ws.Range("D148").Value = "NewEngland"
$ws.Range("E148").Value = 0
$ws.Range("F148").Value = 9
$ws.Range("G148").Value = 6
$ws.Range("H148").Value = 10
$ws.Range("I148").Value = 25
$ws.Range("J148").Value = 7.5
$ws.Range("K148").Value = 9
$ws.Range("L148").Value = -360
$ws.Range("M148").Value = 7

$ws.Range("A149").Value = 1010
$ws.Range("B149").Value = 468
$ws.Range("C149").Value = "H"
$ws.Range("D149").Value = "Houston"
$ws.Range("E149").Value = 6
$ws.Range("F149").Value = 9
$ws.Range("G149").Value = 7
$ws.Range("H149").Value = 0
$ws.Range("I149").Value = 22
$ws.Range("J149").Value = 42
$ws.Range("K149").Value = 39
$ws.Range("L149").Value = 300
$ws.Range("M149").Value = 21

$ws.Range("A150").Value = 1010
$ws.Range("B150").Value = 469
$ws.Range("C150").Value = "V"
$ws.Range("D150").Value = "Chicago"
$ws.Range("E150").Value = 0
$ws.Range("F150").Value = 14
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 6
$ws.Range("I150").Value = 20
$ws.Range("J150").Value = 45.5
$ws.Range("K150").Value = 46
$ws.Range("L150").Value = 200
$ws.Range("M150").Value = 23.5

$ws.Range("A151").Value = 1010
$ws.Range("B151").Value = 470
$ws.Range("C151").Value = "H"
$ws.Range("D151").Value = "LasVegas"
$ws.Range("E151").Value = 3
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 6
$ws.Range("I151").Value = 9
$ws.Range("J151").Value = 7
$ws.Range("K151").Value = 5.5
$ws.Range("L151").Value = -240
$ws.Range("M151").Value = 6.5

$ws.Range("A152").Value = 1010
$ws.Range("B152").Value = 471
$ws.Range("C152").Value = "V"
$ws.Range("D152").Value = "Cleveland"
$ws.Range("E152").Value = 3
$ws.Range("F152").Value = 17
$ws.Range("G152").Value = 7
$ws.Range("H152").Value = 15
$ws.Range("I152").Value = 42
$ws.Range("J152").Value = 1
$ws.Range("K152").Value = 47
$ws.Range("L152").Value = 115
$ws.Range("M152").Value = 26.5

$ws.Range("A153").Value = 1010
$ws.Range("B153").Value = 472
$ws.Range("C153").Value = "H"
$ws.Range("D153").Value = "LAChargers"
$ws.Range("E153").Value = 7
$ws.Range("F153").Value = 6
$ws.Range("G153").Value = 8
$ws.Range("H153").Value = 26
$ws.Range("I153").Value = 47
$ws.Range("J153").Value = 50
$ws.Range("K153").Value = 2.5
$ws.Range("L153").Value = -135
$ws.Range("M153").Value = 3

$ws.Range("A154").Value = 1010
$ws.Range("B154").Value = 473
$ws.Range("C154").Value = "V"
$ws.Range("D154").Value = "NYGiants"
$ws.Range("E154").Value = 0
$ws.Range("F154").Value = 10
$ws.Range("G154").Value = 3
$ws.Range("H154").Value = 7
$ws.Range("I154").Value = 20
$ws.Range("J154").Value = 49
$ws.Range("K154").Value = 53.5
$ws.Range("L154").Value = 290
$ws.Range("M154").Value = 23.5

$ws.Range("A155").Value = 1010
$ws.Range("B155").Value = 474
$ws.Range("C155").Value = "H"
$ws.Range("D155").Value = "Dallas"
$ws.Range("E155").Value = 3
$ws.Range("F155").Value = 14
$ws.Range("G155").Value = 10
$ws.Range("H155").Value = 17
$ws.Range("I155").Value = 44
$ws.Range("J155").Value = 8
$ws.Range("K155").Value = 7.5
$ws.Range("L155").Value = -350
$ws.Range("M155").Value = 5.5

$ws.Range("A156").Value = 1010
$ws.Range("B156").Value = 475
$ws.Range("C156").Value = "V"
$ws.Range("D156").Value = "SanFrancisco"
$ws.Range("E156").Value = 0
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 7
$ws.Range("H156").Value = 3
$ws.Range("I156").Value = 10
$ws.Range("J156").Value = 52
$ws.Range("K156").Value = 48.5
$ws.Range("L156").Value = 235
$ws.Range("M156").Value = 22.5

$ws.Range("A157").Value = 1010
$ws.Range("B157").Value = 476
$ws.Range("C157").Value = "H"
$ws.Range("D157").Value = "Arizona"
$ws.Range("E157").Value = 7
$ws.Range("F157").Value = 3
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 7
$ws.Range("I157").Value = 17
$ws.Range("J157").Value = 2
$ws.Range("K157").Value = 6
$ws.Range("L157").Value = -275
$ws.Range("M157").Value = 1

$ws.Range("A158").Value = 1010
$ws.Range("B158").Value = 477
$ws.Range("C158").Value = "V"
$ws.Range("D158").Value = "Buffalo"
$ws.Range("E158").Value = 7
$ws.Range("F158").Value = 17
$ws.Range("G158").Value = 7
$ws.Range("H158").Value = 7
$ws.Range("I158").Value = 38
$ws.Range("J158").Value = 56.5
$ws.Range("K158").Value = 57
$ws.Range("L158").Value = 120
$ws.Range("M158").Value = 28.5

$ws.Range("A159").Value = 1010
$ws.Range("B159").Value = 478
$ws.Range("C159").Value = "H"
$ws.Range("D159").Value = "KansasCity"
$ws.Range("E159").Value = 3
$ws.Range("F159").Value = 10
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 7
$ws.Range("I159").Value = 20
$ws.Range("J159").Value = 4
$ws.Range("K159").Value = 3
$ws.Range("L159").Value = -140
$ws.Range("M159").Value = 4.5
